$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 95, shifting rows 95:185 down to 96:186
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly record
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C95").Value = 'Ñuble'
$ws.Range("D95").Value = '2022-07-25'
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = 'Fruta'
$ws.Range("G95").Value = 100101
$ws.Range("H95").Value = 'Berries'
$ws.Range("I95").Value = 100101007
$ws.Range("J95").Value = 'Kiwi'
$ws.Range("K95").Value = 'Hayward'
$ws.Range("L95").Value = 'Primera'
$ws.Range("M95").Value = 120
$ws.Range("N95").Value = 6500
$ws.Range("O95").Value = 7000
$ws.Range("P95").Value = 6750
$ws.Range("Q95").Value = '$/bandeja 18 kilos'
$ws.Range("R95").Value = 'Provincia de Curicó'
$ws.Range("S95").Value = 375
$ws.Range("T95").Value = 18
